$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Frutas")

# Insert a new row at the top, shifting existing data down
$ws.Rows.Item(1).Insert()

# Set the new header row values
$ws.Cells.Item(1, 1).Value = "Frutas"
$ws.Cells.Item(1, 2).Value = "Quantidade"
$ws.Cells.Item(1, 3).Value = "Preço"
